$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Helper: writes a literal text value into a cell without Excel's COM
# layer auto-coercing numeric-looking strings (e.g. "66.175.58",
# "28.17", "1.00") into actual numbers. We build a text formula,
# then copy/paste-special as values so the stored cell keeps its
# original "string" cell type and the worksheet's cell styles / the
# workbook's style table stay untouched.
function Set-TextValue($range, $value) {
    $escaped = $value -replace '"', '""'
    $range.Formula = "=""$escaped"""
    $range.Copy()
    $range.PasteSpecial(-4163)
}

Set-TextValue $ws.Range('D2') '65.954.20'
Set-TextValue $ws.Range('E2') '  +1.31%  '
Set-TextValue $ws.Range('D3') '3.587.49'
Set-TextValue $ws.Range('E3') '  +1.43%  '
Set-TextValue $ws.Range('E4') '  +0.04%  '
Set-TextValue $ws.Range('D5') '601.25'
Set-TextValue $ws.Range('E5') '  +0.74%  '
Set-TextValue $ws.Range('D6') '138.79'
Set-TextValue $ws.Range('E6') '  +0.75%  '
Set-TextValue $ws.Range('D7') '3.583.63'
Set-TextValue $ws.Range('E7') '  +1.35%  '
Set-TextValue $ws.Range('E8') '  +0.02%  '
Set-TextValue $ws.Range('D9') '0.499'
Set-TextValue $ws.Range('E9') '  +1.20%  '
Set-TextValue $ws.Range('D10') '0.126'
Set-TextValue $ws.Range('E10') '  +1.46%  '
Set-TextValue $ws.Range('D11') '7.21'
Set-TextValue $ws.Range('E11') '  +4.39%  '
Set-TextValue $ws.Range('D12') '0.393'
Set-TextValue $ws.Range('E12') '  +1.85%  '
Set-TextValue $ws.Range('D13') '4.200.85'
Set-TextValue $ws.Range('E13') '  +1.56%  '
Set-TextValue $ws.Range('D14') '28.17'
Set-TextValue $ws.Range('E14') '  +3.32%  '
Set-TextValue $ws.Range('D15') '0.0000186'
Set-TextValue $ws.Range('E15') '  +1.81%  '
Set-TextValue $ws.Range('D16') '3.594.16'
Set-TextValue $ws.Range('E16') '  +1.43%  '
Set-TextValue $ws.Range('E17') '  +0.35%  '
Set-TextValue $ws.Range('D18') '66.030.43'
Set-TextValue $ws.Range('E18') '  +1.43%  '
Set-TextValue $ws.Range('D19') '10.05'
Set-TextValue $ws.Range('E19') '  +0.08%  '
Set-TextValue $ws.Range('D20') '14.61'
Set-TextValue $ws.Range('E20') '  +2.59%  '
Set-TextValue $ws.Range('D21') '5.85'
Set-TextValue $ws.Range('E21') '  -0.90%  '
Set-TextValue $ws.Range('D22') '397.11'
Set-TextValue $ws.Range('E22') '  +1.18%  '
Set-TextValue $ws.Range('D23') '0.587'
Set-TextValue $ws.Range('E23') '  +2.69%  '
Set-TextValue $ws.Range('D24') '3.734.59'
Set-TextValue $ws.Range('E24') '  +1.50%  '
Set-TextValue $ws.Range('D25') '75.02'
Set-TextValue $ws.Range('E25') '  +2.11%  '
Set-TextValue $ws.Range('E26') '  -0.06%  '
Set-TextValue $ws.Range('D27') '0.0000120'
Set-TextValue $ws.Range('E27') '  +5.62%  '
Set-TextValue $ws.Range('D28') '8.11'
Set-TextValue $ws.Range('E28') '  +3.86%  '
Set-TextValue $ws.Range('E29') '  +24.03%  '
Set-TextValue $ws.Range('D30') '8.64'
Set-TextValue $ws.Range('E30') '  +4.73%  '
Set-TextValue $ws.Range('D31') '2.34'
Set-TextValue $ws.Range('E31') '  +2.83%  '
Set-TextValue $ws.Range('D32') '1.00'
Set-TextValue $ws.Range('E32') '  +0.17%  '
Set-TextValue $ws.Range('D33') '3.594.98'
Set-TextValue $ws.Range('E33') '  +1.02%  '
Set-TextValue $ws.Range('D34') '24.52'
Set-TextValue $ws.Range('E34') '  +3.18%  '
Set-TextValue $ws.Range('E35') '  +3.86%  '
Set-TextValue $ws.Range('E36') '  -0.01%  '
Set-TextValue $ws.Range('D37') '5.41'
Set-TextValue $ws.Range('E37') '  +8.64%  '
Set-TextValue $ws.Range('D38') '1.61'
Set-TextValue $ws.Range('E38') '  +1.65%  '
Set-TextValue $ws.Range('D39') '7.02'
Set-TextValue $ws.Range('E39') '  +1.43%  '
Set-TextValue $ws.Range('D40') '168.82'
Set-TextValue $ws.Range('E40') '  -1.17%  '
Set-TextValue $ws.Range('D41') '0.0838'
Set-TextValue $ws.Range('E41') '  +4.41%  '
Set-TextValue $ws.Range('D42') '0.840'
Set-TextValue $ws.Range('E42') '  +2.19%  '
Set-TextValue $ws.Range('D43') '1.27'
Set-TextValue $ws.Range('E43') '  +6.39%  '
Set-TextValue $ws.Range('D44') '25.97'
Set-TextValue $ws.Range('E44') '  -1.09%  '
Set-TextValue $ws.Range('D45') '43.17'
Set-TextValue $ws.Range('E45') '  +1.45%  '
Set-TextValue $ws.Range('D46') '4.55'
Set-TextValue $ws.Range('E46') '  +2.71%  '
Set-TextValue $ws.Range('E47') '  +0.05%  '
Set-TextValue $ws.Range('D48') '1.71'
Set-TextValue $ws.Range('E48') '  +2.39%  '
Set-TextValue $ws.Range('D49') '6.95'
Set-TextValue $ws.Range('E49') '  +2.05%  '
Set-TextValue $ws.Range('D50') '2.440.70'
Set-TextValue $ws.Range('E50') '  +2.45%  '
Set-TextValue $ws.Range('D51') '0.913'
Set-TextValue $ws.Range('E51') '  +9.13%  '
